$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Link" column in front of the existing table (shifts Firm/Name/Role/
# Country/Nationality/Practice Area/Email/Phone one column to the right).
[void]$ws.Columns("A").Insert()

# New header text for the inserted column.
$ws.Range("A1").Value2 = "Link"

# Match the header formatting (number format, font, fill, border, alignment)
# used by the other "title" style header cells (Firm/Name/Role) by copying it
# from the neighboring header cell that already has it.
[void]$ws.Range("B1").Copy()
[void]$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give the new column roughly the same width as the other title columns.
$ws.Columns("A").ColumnWidth = 26.83

# Move the active selection to A2, like in the edited workbook.
[void]$ws.Range("A2").Select()
